$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "n_estimators: 50/loss: exponential/learning_rate: 0.3"
$ws.Range("D2").Value = 0.9968105280705966
$ws.Range("E2").Value = 0.9957473707607954
$ws.Range("F2").Value = 6.109256117647064
$ws.Range("G2").Value = 2.488673687089965
$ws.Range("H2").Value = 1.577553069500347
$ws.Range("I2").Value = 8.364261831033858
$ws.Range("J2").Value = 2.892103357598732
$ws.Range("K2").Value = 0.2155627060651892

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "n_estimators: 100/max_samples: 1.0"
$ws.Range("D3").Value = 0.9988560555848621
$ws.Range("E3").Value = 0.9984747407798161
$ws.Range("F3").Value = 3.425956200000073
$ws.Range("G3").Value = 1.322533833529427
$ws.Range("H3").Value = 1.150014710136104
$ws.Range("I3").Value = 2.999948210910336
$ws.Range("J3").Value = 1.732035857281926
$ws.Range("K3").Value = 0.1316424817931592

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "default"
$ws.Range("D4").Value = 0.9862461086791027
$ws.Range("E4").Value = 0.9816614782388036
$ws.Range("F4").Value = 9.578972520548049
$ws.Range("G4").Value = 5.554871476977115
$ws.Range("H4").Value = 2.356877484507227
$ws.Range("I4").Value = 36.06902670721665
$ws.Range("J4").Value = 6.005749470900085
$ws.Range("K4").Value = 0.4012172529941787

# Row 5
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 0.9979738033433034
$ws.Range("E5").Value = 0.9972984044577379
$ws.Range("F5").Value = 4.138573999999991
$ws.Range("G5").Value = 1.993801235294116
$ws.Range("H5").Value = 1.412020267309969
$ws.Range("I5").Value = 5.313619223777232
$ws.Range("J5").Value = 2.305128895263176
$ws.Range("K5").Value = 0.1371398890472471

# Row 6
$ws.Range("B6").Value = 2
$ws.Range("D6").Value = -0.00005856187952280045
$ws.Range("E6").Value = -0.3334114158393637
$ws.Range("F6").Value = 79.26752778125
$ws.Range("G6").Value = 41.21908382904412
$ws.Range("H6").Value = 6.420209017551073
$ws.Range("I6").Value = 2510.325650115235
$ws.Range("J6").Value = 50.10315010171751
$ws.Range("K6").Value = 3.827869832326928

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "normalize: False/l1_ratio: 0.75/fit_intercept: False/alpha: 0.02"
$ws.Range("D7").Value = 0.9542793937512544
$ws.Range("E7").Value = 0.9390391916683392
$ws.Range("F7").Value = 19.20629528732702
$ws.Range("G7").Value = 9.006199801404126
$ws.Range("H7").Value = 3.001033122343725
$ws.Range("I7").Value = 119.9004506710438
$ws.Range("J7").Value = 10.94990642293549
$ws.Range("K7").Value = 0.5675978780939747

# Row 8
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "default"
$ws.Range("D8").Value = 0.9842965028627755
$ws.Range("E8").Value = 0.979062003817034
$ws.Range("F8").Value = 9.819193363552763
$ws.Range("G8").Value = 5.825573130736727
$ws.Range("H8").Value = 2.413622408484129
$ws.Range("I8").Value = 39.41858323376366
$ws.Range("J8").Value = 6.278422033740934
$ws.Range("K8").Value = 2.555655046335292

# Row 9
$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 0.9986523022906715
$ws.Range("E9").Value = 0.9982030697208953
$ws.Range("G9").Value = 1.411778411764706
$ws.Range("H9").Value = 1.188182819167449
$ws.Range("I9").Value = 3.534283028481468
$ws.Range("J9").Value = 1.879968890296185
$ws.Range("K9").Value = 0.04215737014596937

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 0.9992304978263927
$ws.Range("E10").Value = 0.998973997101857
$ws.Range("F10").Value = 2.795443600000084
$ws.Range("G10").Value = 1.158941678823537
$ws.Range("H10").Value = 1.076541536041939
$ws.Range("I10").Value = 2.017988495294463
$ws.Range("J10").Value = 1.420559219214202
$ws.Range("K10").Value = 0.06739951589098545

# Row 11
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 0.9989388349822468
$ws.Range("E11").Value = 0.9985851133096624
$ws.Range("F11").Value = 5.394667161180289
$ws.Range("G11").Value = 1.132850583870979
$ws.Range("H11").Value = 1.064354538615296
$ws.Range("I11").Value = 2.782862571260226
$ws.Range("J11").Value = 1.668191407261237
$ws.Range("K11").Value = 0.2001410298244644

# Row 12
$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 0.9855508370577364
$ws.Range("E12").Value = 0.9807344494103152
$ws.Range("F12").Value = 10.80352007121122
$ws.Range("G12").Value = 5.440576472896923
$ws.Range("H12").Value = 2.332504335022108
$ws.Range("I12").Value = 36.26998031844126
$ws.Range("J12").Value = 6.022456335951408
$ws.Range("K12").Value = 2.759056682813132

# Row 13
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "default"
$ws.Range("D13").Value = 0.9337517906779914
$ws.Range("E13").Value = 0.9116690542373219
$ws.Range("F13").Value = 27.6298188
$ws.Range("G13").Value = 10.96940407058823
$ws.Range("H13").Value = 3.312009068615036
$ws.Range("I13").Value = 173.7332639607421
$ws.Range("J13").Value = 13.18079147702224
$ws.Range("K13").Value = 0.870974381660599

# Row 14
$ws.Range("B14").Value = 1
$ws.Range("D14").Value = 0.9862975178373595
$ws.Range("E14").Value = 0.981730023783146
$ws.Range("F14").Value = 9.527663346272767
$ws.Range("G14").Value = 5.547959643245941
$ws.Range("H14").Value = 2.35541071646665
$ws.Range("I14").Value = 35.93420825773965
$ws.Range("J14").Value = 5.994514847570206
$ws.Range("K14").Value = 0.4012007676083679

# Row 15
$ws.Range("B15").Value = 2
$ws.Range("D15").Value = 0.9859272148331321
$ws.Range("E15").Value = 0.981236286444176
$ws.Range("F15").Value = 10.82754213475699
$ws.Range("G15").Value = 5.419625058810348
$ws.Range("H15").Value = 2.328008818456311
$ws.Range("I15").Value = 35.3252048625585
$ws.Range("J15").Value = 5.943501061037888
$ws.Range("K15").Value = 2.477894955344792

# Row 16
$ws.Range("B16").Value = 2
$ws.Range("D16").Value = 0.9863099920604954
$ws.Range("E16").Value = 0.9817466560806606
$ws.Range("F16").Value = 10.82413047261112
$ws.Range("G16").Value = 5.396435320041327
$ws.Range("H16").Value = 2.323022884097642
$ws.Range("I16").Value = 34.36436563897902
$ws.Range("J16").Value = 5.86211272827289
$ws.Range("K16").Value = 2.094464284029546

# Row 17
$ws.Range("B17").Value = 2
$ws.Range("D17").Value = 0.9858865666226108
$ws.Range("E17").Value = 0.9811820888301477
$ws.Range("F17").Value = 10.81247518970051
$ws.Range("G17").Value = 5.419353069239344
$ws.Range("H17").Value = 2.327950400940566
$ws.Range("I17").Value = 35.42723913274293
$ws.Range("J17").Value = 5.952078555659605
$ws.Range("K17").Value = 2.506906523463603

# Row 18
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "default"
$ws.Range("D18").Value = 0.9604634768935373
$ws.Range("E18").Value = 0.9472846358580497
$ws.Range("F18").Value = 18.00584804673306
$ws.Range("G18").Value = 8.429227032204874
$ws.Range("H18").Value = 2.903313113015004
$ws.Range("I18").Value = 99.24373617080612
$ws.Range("J18").Value = 9.962115045049726
$ws.Range("K18").Value = 2.057557845008537

# Row 19
$ws.Range("B19").Value = 2
$ws.Range("D19").Value = 0.9859272148331321
$ws.Range("E19").Value = 0.981236286444176
$ws.Range("F19").Value = 10.82754213475699
$ws.Range("G19").Value = 5.419625058810348
$ws.Range("H19").Value = 2.328008818456311
$ws.Range("I19").Value = 35.3252048625585
$ws.Range("J19").Value = 5.943501061037888
$ws.Range("K19").Value = 2.477894955344792

# Row 20
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "normalize: False/fit_intercept: True/criterion: bic"
$ws.Range("D20").Value = 0.9855178157563241
$ws.Range("E20").Value = 0.9806904210084322
$ws.Range("F20").Value = 10.80788753395722
$ws.Range("G20").Value = 5.449274755842302
$ws.Range("H20").Value = 2.334368170585416
$ws.Range("I20").Value = 36.35286968421935
$ws.Range("J20").Value = 6.029334099568488
$ws.Range("K20").Value = 2.669498630822106

# Row 21
$ws.Range("B21").Value = 1
$ws.Range("D21").Value = 0.9862975178373594
$ws.Range("E21").Value = 0.9817300237831459
$ws.Range("F21").Value = 9.527663346272874
$ws.Range("G21").Value = 5.547959643245955
$ws.Range("H21").Value = 2.355410716466654
$ws.Range("I21").Value = 35.93420825773995
$ws.Range("J21").Value = 5.994514847570231
$ws.Range("K21").Value = 0.4012007676083682

# Row 22
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "loss: squared_epsilon_insensitive/epsilon: 0/C: 5"
$ws.Range("D22").Value = 0.98552249644859
$ws.Range("E22").Value = 0.9806966619314533
$ws.Range("F22").Value = 10.35094127830876
$ws.Range("G22").Value = 5.47579981558126
$ws.Range("H22").Value = 2.340042695247516
$ws.Range("I22").Value = 36.34112031042964
$ws.Range("J22").Value = 6.028359669962438
$ws.Range("K22").Value = 2.797194473915273

# Row 23
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "nu: 0.01/kernel: sigmoid/gamma: 0.0001/C: 100000"
$ws.Range("D23").Value = 0.1917950719520626
$ws.Range("E23").Value = -0.0776065707305833
$ws.Range("F23").Value = 68.05770179866666
$ws.Range("G23").Value = 35.75978546768628
$ws.Range("H23").Value = 5.979948617478772
$ws.Range("I23").Value = 2028.738754673741
$ws.Range("J23").Value = 45.04152256167348
$ws.Range("K23").Value = 1.25712151558167

# Row 24
$ws.Range("B24").Value = 2
$ws.Range("D24").Value = 0.9471168534958073
$ws.Range("E24").Value = 0.9294891379944098
$ws.Range("F24").Value = 20.1731404857955
$ws.Range("G24").Value = 10.18639173713235
$ws.Range("H24").Value = 3.191612717284532
$ws.Range("I24").Value = 132.7461452645105
$ws.Range("J24").Value = 11.5215513393167
$ws.Range("K24").Value = 3.343254092502247

# Row 25
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = 0.9861505789527807
$ws.Range("E25").Value = 0.9815341052703742
$ws.Range("F25").Value = 9.555304033238144
$ws.Range("G25").Value = 5.634958671136338
$ws.Range("H25").Value = 2.373806788922876
$ws.Range("I25").Value = 36.31954957159367
$ws.Range("J25").Value = 6.026570299232696
$ws.Range("K25").Value = 0.408824060076549

# Row 26
$ws.Range("B26").Value = 2
$ws.Range("D26").Value = 0.9852282362272531
$ws.Range("E26").Value = 0.9803043149696707
$ws.Range("F26").Value = 13.74633211974675
$ws.Range("G26").Value = 5.087720810867804
$ws.Range("H26").Value = 2.255597661567285
$ws.Range("I26").Value = 37.07976603537834
$ws.Range("J26").Value = 6.089315728009046
$ws.Range("K26").Value = 1.885423497576516

# Row 27
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "min_samples: 1/loss: absolute_loss"
$ws.Range("D27").Value = 0.9868373379031185
$ws.Range("E27").Value = 0.9824497838708246
$ws.Range("F27").Value = 10.22669917735227
$ws.Range("G27").Value = 5.184240702942263
$ws.Range("H27").Value = 2.27689277370329
$ws.Range("I27").Value = 33.04063336401102
$ws.Range("J27").Value = 5.74809823889702
$ws.Range("K27").Value = 1.747650457449339

# Row 28
$ws.Range("B28").Value = 1
$ws.Range("D28").Value = 0.9988428424442749
$ws.Range("E28").Value = 0.9984571232590332
$ws.Range("F28").Value = 3.630354080000075
$ws.Range("G28").Value = 1.334194305882373
$ws.Range("H28").Value = 1.155073290264463
$ws.Range("I28").Value = 3.03459914057142
$ws.Range("J28").Value = 1.742010086242735
$ws.Range("K28").Value = 0.1427371260241757

# Row 29
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = "default"
$ws.Range("D29").Value = 0.9860761163890228
$ws.Range("E29").Value = 0.9814348218520305
$ws.Range("F29").Value = 8.868482096518939
$ws.Range("G29").Value = 5.573437117941317
$ws.Range("H29").Value = 2.360812808746453
$ws.Range("I29").Value = 34.95143535610881
$ws.Range("J29").Value = 5.911973896771602
$ws.Range("K29").Value = 2.72801337892266

# Row 30
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = "normalize: False/gcv_mode: eigen/fit_intercept: True"
$ws.Range("D30").Value = 0.9855386902963195
$ws.Range("E30").Value = 0.9807182537284259
$ws.Range("F30").Value = 10.36384665501048
$ws.Range("G30").Value = 5.47199638603862
$ws.Range("H30").Value = 2.339229870286077
$ws.Range("I30").Value = 36.30047085960842
$ws.Range("J30").Value = 6.02498720825268
$ws.Range("K30").Value = 2.793948536557211

# Row 31
$ws.Range("B31").Value = 1
$ws.Range("D31").Value = 0.9862020246423787
$ws.Range("E31").Value = 0.9816026995231716
$ws.Range("F31").Value = 9.651984400222851
$ws.Range("G31").Value = 5.561411716293115
$ws.Range("H31").Value = 2.358264556043939
$ws.Range("I31").Value = 36.18463532014376
$ws.Range("J31").Value = 6.015366598981625
$ws.Range("K31").Value = 0.4021119255477097

# Row 32
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "kernel: sigmoid/gamma: 0.01/epsilon: 1.0/C: 5000"
$ws.Range("D32").Value = 0.9837955433770645
$ws.Range("E32").Value = 0.978394057836086
$ws.Range("F32").Value = 10.5882702982145
$ws.Range("G32").Value = 5.523941464245285
$ws.Range("H32").Value = 2.350306674509794
$ws.Range("I32").Value = 42.49553563219968
$ws.Range("J32").Value = 6.51885999483036
$ws.Range("K32").Value = 0.4787901291670473

# Row 33
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "learning_rate: 0.05/booster: gbtree"
$ws.Range("D33").Value = 0.9987470923905283
$ws.Range("E33").Value = 0.9983294565207044
$ws.Range("F33").Value = 3.628339296630859
$ws.Range("G33").Value = 1.307016884193981
$ws.Range("H33").Value = 1.143248391292977
$ws.Range("I33").Value = 3.28569980475623
$ws.Range("J33").Value = 1.812649939937723
$ws.Range("K33").Value = 0.08495895651026102
